# Insert a new weekly price record as row 42 on the daily-logic sheet
# (Fruta, Terminal La Palmera de La Serena - Caqui). Everything currently
# at row 42 and below shifts down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("42").Insert()

$ws.Range("A42").Value = 8
$ws.Range("B42").Value = "Terminal La Palmera de La Serena"
$ws.Range("C42").Value = "Coquimbo"
$ws.Range("D42").Value = 45071
$ws.Range("E42").Value = 4
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100107
$ws.Range("H42").Value = "Otros"
$ws.Range("I42").Value = 100107001
$ws.Range("J42").Value = "Caqui"
$ws.Range("K42").Value = "Mankaki"
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 16
$ws.Range("N42").Value = 310000
$ws.Range("O42").Value = 320000
$ws.Range("P42").Value = 315000
$ws.Range("Q42").Value = "`$/bins (450 kilos)"
$ws.Range("R42").Value = "Región de O'Higgins"
$ws.Range("S42").Value = 700
$ws.Range("T42").Value = 450
